{"js": "// Insert a new bulleted list item right after the paragraph that ends with\n// \"...and if an empty line was given, it will be ignored.\" and before the\n// paragraph that begins with \"If the level is completely empty...\".\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst anchorText = \"and if an empty line was given, it will be ignored.\";\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(anchorText) !== -1) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Could not find the anchor paragraph ending in: \" + anchorText);\n}\n\n// InsertParagraphAfter inherits the anchor paragraph's list/paragraph\n// formatting (ListParagraph style, numId 2 / ilvl 0) and run formatting\n// (en-US language), matching the target markup exactly.\nanchor.insertParagraph(\n  \"If a character is wrong or a space was given, the character will be ignored.\",\n  \"After\"\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the paragraph that ends with \"...and if an empty line was given,\n# it will be ignored.\" by scanning the document's paragraphs collection.\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"*and if an empty line was given, it will be ignored.*\") {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Could not find the anchor paragraph ending in: and if an empty line was given, it will be ignored.\"\n}\n\n$anchor = $d.Paragraphs.Item($anchorIndex)\n\n# InsertParagraphAfter copies the anchor paragraph's formatting (ListParagraph\n# style, numId 2 / ilvl 0 list membership, en-US run language) onto the new,\n# empty paragraph - matching the target markup. Then grab that freshly minted\n# paragraph (immediately after the anchor) and fill in its text.\n$anchor.Range.InsertParagraphAfter()\n$newPara = $d.Paragraphs.Item($anchorIndex + 1)\n$newPara.Range.Text = \"If a character is wrong or a space was given, the character will be ignored.\"\n"}
